# Update "想去人数" (column F) values on several rows across the
# "展览", "演出", "本地生活" and "全部类型" worksheets to reflect a newer
# scrape of the source data (commit: "Update gh-pages to output generated
# at 456a3b4").
#
# Note: this runtime's PowerShell-style function parameter binding only
# works reliably with *positional* arguments, so Set-FValue is called
# positionally below (named args like -SheetName bind empty strings).

$wb = $excel.ActiveWorkbook

function Set-FValue {
    param(
        [string]$SheetName,
        [string]$CellRef,
        [double]$NewValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $NewValue
}

# Sheet "展览"
Set-FValue "展览" "F3"  963
Set-FValue "展览" "F7"  1181
Set-FValue "展览" "F9"  36
Set-FValue "展览" "F11" 1049
Set-FValue "展览" "F12" 2497
Set-FValue "展览" "F15" 1680
Set-FValue "展览" "F23" 770
Set-FValue "展览" "F24" 651
Set-FValue "展览" "F25" 514
Set-FValue "展览" "F28" 35
Set-FValue "展览" "F31" 326
Set-FValue "展览" "F34" 1414
Set-FValue "展览" "F35" 468
Set-FValue "展览" "F38" 4029

# Sheet "演出"
Set-FValue "演出" "F31" 1726
Set-FValue "演出" "F37" 5

# Sheet "本地生活"
Set-FValue "本地生活" "F5" 1681
Set-FValue "本地生活" "F7" 1031

# Sheet "全部类型"
Set-FValue "全部类型" "F3"  1681
Set-FValue "全部类型" "F7"  963
Set-FValue "全部类型" "F9"  1181
Set-FValue "全部类型" "F12" 36
Set-FValue "全部类型" "F17" 1049
Set-FValue "全部类型" "F19" 2497
Set-FValue "全部类型" "F21" 1680
Set-FValue "全部类型" "F28" 770
Set-FValue "全部类型" "F29" 651
Set-FValue "全部类型" "F30" 514
Set-FValue "全部类型" "F32" 35
Set-FValue "全部类型" "F37" 326
Set-FValue "全部类型" "F44" 1726
Set-FValue "全部类型" "F45" 1726
Set-FValue "全部类型" "F46" 1414
Set-FValue "全部类型" "F49" 4029
Set-FValue "全部类型" "F50" 5
